$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data row (C2 first so shared-string order matches)
$ws.Range("C2").Value = "AdminSaurabTest@gmail.com"

# Update header row
$ws.Range("C1").Value = "emailid"
$ws.Range("D1").Value = "days"
$ws.Range("E1").Value = "month "

# Remaining data row values
$ws.Range("D2").Value = 20
$ws.Range("E2").Value = 8

# Widen column C to fit the new "emailid" data (target ~28.11 chars).
# This runtime snaps ColumnWidth to the nearest 1/6-character increment,
# so 27.3 is the input that resolves to the closest achievable stored
# width (28.1666...) to the authored value of 28.109375.
$ws.Range("C1").ColumnWidth = 27.3

# Update selection to F6
$ws.Range("F6").Select()
